# Map135.xlsx "huge v2 update" patch:
# Duplicate column A (translation strings) into a new column B,
# row for row, for rows 1 through 98, reusing the same shared-string
# values (so no new unique strings get added to the shared string table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 98

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Text
}
